# Fruta / hortaliza, semanal
# Rows 2-15 hold one weekly price record each (Femacal de La Calera - Frambuesa).
# This refresh reshuffles the per-row Fecha/Volumen/Precio.../Origen values
# (rows 4 and 8 are unchanged) to match the updated weekly extract.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44582
$ws.Range("M2").Value = 150
$ws.Range("N2").Value = 6000
$ws.Range("O2").Value = 6500
$ws.Range("P2").Value = 6233
$ws.Range("R2").Value = "Provincia de Curicó"
$ws.Range("S2").Value = 3116

$ws.Range("D3").Value = 44627
$ws.Range("M3").Value = 45
$ws.Range("N3").Value = 6000
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 6000
$ws.Range("R3").Value = "Provincia de Linares"
$ws.Range("S3").Value = 3000

$ws.Range("D5").Value = 44588
$ws.Range("M5").Value = 160
$ws.Range("N5").Value = 6500
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 6750
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 3375

$ws.Range("D6").Value = 44606
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 7000
$ws.Range("O6").Value = 7000
$ws.Range("P6").Value = 7000
$ws.Range("R6").Value = "Provincia de Linares"
$ws.Range("S6").Value = 3500

$ws.Range("D7").Value = 44587
$ws.Range("M7").Value = 165
$ws.Range("N7").Value = 6500
$ws.Range("O7").Value = 7000
$ws.Range("P7").Value = 6742
$ws.Range("R7").Value = "Provincia de Linares"
$ws.Range("S7").Value = 3371

$ws.Range("D9").Value = 44214
$ws.Range("M9").Value = 48
$ws.Range("N9").Value = 6000
$ws.Range("O9").Value = 6000
$ws.Range("P9").Value = 6000
$ws.Range("R9").Value = "Provincia de Linares"
$ws.Range("S9").Value = 3000

$ws.Range("D10").Value = 44614
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 6000
$ws.Range("O10").Value = 6000
$ws.Range("P10").Value = 6000
$ws.Range("R10").Value = "Provincia de Linares"
$ws.Range("S10").Value = 3000

$ws.Range("D11").Value = 44209
$ws.Range("M11").Value = 58
$ws.Range("N11").Value = 6000
$ws.Range("O11").Value = 6000
$ws.Range("P11").Value = 6000
$ws.Range("R11").Value = "Provincia de Curicó"
$ws.Range("S11").Value = 3000

$ws.Range("D12").Value = 44585
$ws.Range("M12").Value = 160
$ws.Range("N12").Value = 6500
$ws.Range("O12").Value = 7000
$ws.Range("P12").Value = 6750
$ws.Range("R12").Value = "Provincia de Curicó"
$ws.Range("S12").Value = 3375

$ws.Range("D13").Value = 44211
$ws.Range("M13").Value = 45
$ws.Range("N13").Value = 6000
$ws.Range("O13").Value = 6000
$ws.Range("P13").Value = 6000
$ws.Range("R13").Value = "Provincia de Curicó"
$ws.Range("S13").Value = 3000

$ws.Range("D14").Value = 44592
$ws.Range("M14").Value = 30
$ws.Range("N14").Value = 8000
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 8000
$ws.Range("R14").Value = "Provincia de Linares"
$ws.Range("S14").Value = 4000

$ws.Range("D15").Value = 44589
$ws.Range("M15").Value = 60
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 6000
$ws.Range("P15").Value = 6000
$ws.Range("R15").Value = "Provincia de Curicó"
$ws.Range("S15").Value = 3000
